$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7
$ws.Range("F3").Value = 268
$ws.Range("F4").Value = 140
$ws.Range("F7").Value = 4879
$ws.Range("F8").Value = 4879
$ws.Range("E11").Value = "2024.06.29 10:00-06.29 16:00"
$ws.Range("F11").Value = 466
$ws.Range("F13").Value = 1111
$ws.Range("F14").Value = 653
$ws.Range("F15").Value = 4573
$ws.Range("F16").Value = 184
$ws.Range("F17").Value = 191
$ws.Range("F18").Value = 84
$ws.Range("F20").Value = 3588
$ws.Range("F23").Value = 29
$ws.Range("F24").Value = 3331
$ws.Range("F25").Value = 151
$ws.Range("F26").Value = 140
$ws.Range("F28").Value = 348
$ws.Range("F29").Value = 164
$ws.Range("F30").Value = 211
$ws.Range("F33").Value = 74
$ws.Range("F37").Value = 5822
$ws.Range("F38").Value = 913
$ws.Range("F39").Value = 430
$ws.Range("F42").Value = 56
$ws.Range("F43").Value = 1188
$ws.Range("F44").Value = 539
$ws.Range("F45").Value = 20
$ws.Range("F46").Value = 2059
$ws.Range("F49").Value = 729
$ws.Range("F50").Value = 872

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 15
$ws.Range("F5").Value = 2
$ws.Range("F9").Value = 44
$ws.Range("F10").Value = 69
$ws.Range("F12").Value = 63
$ws.Range("F16").Value = 128
$ws.Range("F24").Value = 763

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 15
$ws.Range("F4").Value = 268
$ws.Range("F5").Value = 2
$ws.Range("F7").Value = 140
$ws.Range("F10").Value = 4879
$ws.Range("F11").Value = 4879
$ws.Range("F13").Value = 44
$ws.Range("F14").Value = 69
$ws.Range("E15").Value = "2024.06.29 10:00-06.29 16:00"
$ws.Range("F15").Value = 466
$ws.Range("F16").Value = 1111
$ws.Range("F17").Value = 653
$ws.Range("F18").Value = 4573
$ws.Range("F19").Value = 184
$ws.Range("F20").Value = 191
$ws.Range("F21").Value = 84
$ws.Range("F23").Value = 3588
$ws.Range("F24").Value = 3331
$ws.Range("F25").Value = 151
$ws.Range("F26").Value = 140
$ws.Range("F27").Value = 211
$ws.Range("F30").Value = 74
$ws.Range("F33").Value = 128
$ws.Range("F35").Value = 5822
$ws.Range("F36").Value = 913
$ws.Range("F37").Value = 430
$ws.Range("F42").Value = 56
$ws.Range("F43").Value = 1188
$ws.Range("F44").Value = 539
$ws.Range("F45").Value = 2059
$ws.Range("F48").Value = 729
$ws.Range("F49").Value = 872
